# Predicting Fingers - 2 b
#
# The "constraints" bullet ("What are the constraints? Not double counting
# on your pinky.") loses the _GoBack bookmark that used to wrap its answer,
# and the following "sub-goals" bullet gains a new blue answer ("Figuring
# out if you land on the same finger counting to 100 and 1000 as you do
# counting to 10.") which is now wrapped by that _GoBack bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Drop the old bookmark around "Not double counting on your pinky."
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# 2. Find the "What are the sub-goals?" paragraph that immediately
#    follows the constraints bullet containing "Not double counting on
#    your pinky." (robust to absolute paragraph-index drift).
# ---------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "^What are the sub-goals\?\r?$") {
        $prev = $p.Previous()
        if ($prev -ne $null -and $prev.Range.Text -match "Not double counting on your pinky") {
            $targetPara = $p
        }
    }
}

# ---------------------------------------------------------------------
# 3. Build the new answer text in a scratch paragraph directly after the
#    target bullet, then fold it back in by deleting the intervening
#    paragraph mark. Doing the typing in a fresh paragraph (rather than
#    appending in place) keeps "What are the sub-goals?", the following
#    space, and the new colored answer as three distinct runs -- exactly
#    like the separate runs already used for its sibling bullets --
#    instead of Word's usual same-formatting typing merge.
# ---------------------------------------------------------------------
$r = $targetPara.Range
$endPoint = $d.Range($r.End - 1, $r.End - 1)
$endPoint.InsertParagraphAfter()

$newParaIndex = $targetPara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)
$rNew = $newPara.Range
$insertPt = $d.Range($rNew.Start, $rNew.Start)
$insertPt.InsertAfter(" ")

$rNew2 = $d.Paragraphs($newParaIndex).Range
$coloredStart = $rNew2.End - 1
$insertPt2 = $d.Range($coloredStart, $coloredStart)
$insertPt2.InsertAfter("Figuring out if you land on the same finger counting to 100 and 1000 as you do counting to 10.")

$rNew3 = $d.Paragraphs($newParaIndex).Range
$coloredEnd = $rNew3.End - 1
$coloredRange = $d.Range($coloredStart, $coloredEnd)
$coloredRange.Font.Color = 16737843

# Re-create _GoBack around the new answer text.
$d.Bookmarks.Add("_GoBack", $coloredRange)

# ---------------------------------------------------------------------
# 4. Merge the scratch paragraph back into the sub-goals bullet by
#    deleting the paragraph mark that separates them.
# ---------------------------------------------------------------------
$rTarget = $targetPara.Range
$markRange = $d.Range($rTarget.End - 1, $rTarget.End)
$markRange.Delete()
